$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ink / Varnish rows (8-10): rotate the ink colour assigned to each row.
# Row 11 (Magenta) stays the same.
$ws.Cells.Item(8, 4).Value = "Yellow - Sheet-fed Offset UV - "
$ws.Cells.Item(9, 4).Value = "Cyan - Sheet-fed Offset UV - "
$ws.Cells.Item(10, 4).Value = "Black - Sheet-fed Offset UV - "

# Pallet rows (12-13): swap which job ("Element") each pallet row refers to,
# along with its quantity. The Quantity column stores numeric-looking text
# (e.g. "2.00") as a text string in the source file, so force text format
# before assigning to avoid Excel auto-converting it to a number.
$ws.Cells.Item(12, 2).Value = "1) 591345 2p Packed"
$ws.Cells.Item(12, 5).NumberFormat = "@"
$ws.Cells.Item(12, 5).Value = "2.00"
$ws.Cells.Item(13, 2).Value = "2) 591346 2p Packed"
$ws.Cells.Item(13, 5).NumberFormat = "@"
$ws.Cells.Item(13, 5).Value = "1.00"
